$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.735.22"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.896.78"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.31"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4582"
$ws.Range("E7").Value = "  -1.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3874"
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.72"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07891"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("E11").Value = "  +3.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.78"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Value = "1.907.43"
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.020"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.714"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06956"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.29"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001002"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.03"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "28.767.48"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.311"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "2.113.27"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.061"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.88"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.32"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.808"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.49"
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.914"
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09315"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9244"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.323"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.336"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.257"
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05749"
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.162"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02065"
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.734"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5635"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1786"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.752"
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.210"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07170"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.71"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5321"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.120"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.831"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.49"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.460"
$ws.Range("E51").Value = "  +4.65%  "
